# Applies the recorded edits to A787.xlsx:
#  - Overhead sheet: split "CO LO - gan voi screw" text into separate
#    "screw" rows (12/13/14) and simplified "CO LO (...)" labels (rows 4-7),
#    row-height cleanup on rows 4-7 (back to default) and row 14 (new custom height)
#  - Pax seat sheet: fill in a large block of new part-number rows (4-13),
#    re-tag some A/C rows from A321 -> ATR72, append 6 new blank A321 rows
#    (18-23), widen column E, adjust a couple of row heights
#  - Workbook: Pax seat becomes the active/selected tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overhead sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overhead")

# Rows 4-7: "CO LO - GAN VOI SCREW PN BACS12ER3K7 (...)" -> "CO LO (...)"
$ws1.Range("E4").Value = "CÓ LỖ (24INCH - 35LB VỚI CTR, 50LB VỚI OUTB)"
$ws1.Range("E5").Value = "CÓ LỖ (36INCH - 56LB VỚI CTR, 80LB VỚI OUTB)"
$ws1.Range("E6").Value = "CÓ LỖ (42INCH - 66LB VỚI CTR, 95LB VỚI OUTB)"
$ws1.Range("E7").Value = "CÓ LỖ (48INCH - 77LB VỚI CTR, 111LB VỚI OUTB)"

# Rows 4-7 lose their custom (wrapped-text) row height, back to the sheet default
$ws1.Rows(4).AutoFit()
$ws1.Rows(5).AutoFit()
$ws1.Rows(6).AutoFit()
$ws1.Rows(7).AutoFit()

# New "SCREW" / "KEO DAN TRIM" / "SCREW - TORQUE TUBE" rows
$ws1.Range("B12").Value = "SCREW - TRIM CÓ LỖ"
$ws1.Range("C12").Value = "BACS12ER3K7"

$ws1.Range("B13").Value = "KEO DÁN TRIM"
$ws1.Range("C13").Value = "EC3532B/A-40Z"

$ws1.Range("B14").Value = "SCREW - TORQUE TUBE "
$ws1.Range("C14").Value = "BACS12FA3K3"
$ws1.Range("D14").Value = "BACS12JL3K3"
$ws1.Rows(14).RowHeight = 30

# Selection moves from E16 to E17 on this (now inactive) sheet
$ws1.Range("E17").Select()

# ---------------------------------------------------------------------
# Pax seat sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Pax seat")

# Column E gets its own (wider) width, split off from the old D:E group
$ws3.Columns(5).ColumnWidth = 33.7109375

# Row 4 - ESCUTCHEON
$ws3.Range("B4").Value = "Y SEAT"
$ws3.Range("C4").Value = "ESCUTCHEON"
$ws3.Range("D4").Value = "1010206-301BHU"
$ws3.Range("E4").Value = "1010206-301BHU-VAE / 1010206-301BHUA"
$ws3.Rows(4).RowHeight = 28.5

# Row 5 - ARMCAP GHE IAT (RH)
$ws3.Range("B5").Value = "Y SEAT"
$ws3.Range("C5").Value = "ARMCAP GHẾ IAT"
$ws3.Range("D5").Value = "1012042-071KA02"
$ws3.Range("F5").Value = "HÀNG GHẾ ABCEF (CÁNH TAY MỞ RH)"

# Row 6 - ARMCAP GHE IAT (LH)
$ws3.Range("B6").Value = "Y SEAT"
$ws3.Range("C6").Value = "ARMCAP GHẾ IAT"
$ws3.Range("D6").Value = "1012042-072KA02"
$ws3.Range("F6").Value = "HÀNG GHẾ GHKD (CÁNH TAY MỞ LH)"

# Row 7 - ULTRALOC
$ws3.Range("B7").Value = "Y SEAT"
$ws3.Range("C7").Value = "ULTRALOC"
$ws3.Range("D7").Value = "SP23824L"
$ws3.Range("E7").Value = "1011535-005"

# Row 8 - NUT BAM RECLINE (GHE IAT) -- also A8 goes ATR72 -> B787
$ws3.Range("A8").Value = "B787"
$ws3.Range("B8").Value = "Y SEAT"
$ws3.Range("C8").Value = "NÚT BẤM RECLINE"
$ws3.Range("D8").Value = "SP23997C29"
$ws3.Range("F8").Value = "GHẾ IAT"

# Row 9 - NUT BAM RECLINE (GHE THUONG) -- also A9 goes ATR72 -> B787
$ws3.Range("A9").Value = "B787"
$ws3.Range("B9").Value = "Y SEAT"
$ws3.Range("C9").Value = "NÚT BẤM RECLINE"
$ws3.Range("D9").Value = "SP23997C30"
$ws3.Range("F9").Value = "GHẾ THƯỜNG"

# Row 10 - SEATBELT (IAT hang 16) -- also A10 goes ATR72 -> B787
$ws3.Range("A10").Value = "B787"
$ws3.Range("B10").Value = "Y SEAT"
$ws3.Range("C10").Value = "SEATBELT"
$ws3.Range("D10").Value = "2185-1-052-8022"
$ws3.Range("F10").Value = "GHẾ IAT HÀNG 16 ABC VÀ 16 DEF (TRỪ TÀU 787-10)"
$ws3.Rows(10).RowHeight = 28.5

# Row 11 - SEATBELT (IAT 16 GHK) -- also A11 goes ATR72 -> B787
$ws3.Range("A11").Value = "B787"
$ws3.Range("B11").Value = "Y SEAT"
$ws3.Range("C11").Value = "SEATBELT"
$ws3.Range("D11").Value = "2185-1-062-8022"
$ws3.Range("F11").Value = "GHẾ IAT 16 GHK (TRỪ TÀU 787-10)"

# Row 12 - SEATBELT (Y con lai) -- also A12 goes A321 -> B787, row becomes taller too
$ws3.Range("A12").Value = "B787"
$ws3.Range("B12").Value = "Y SEAT"
$ws3.Range("C12").Value = "SEATBELT"
$ws3.Range("D12").Value = "2006-1-511-8022"
$ws3.Range("F12").Value = "GHẾ Y CÒN LẠI (APPLY CHO TOÀN BỘ GHẾ Y 787-10)"
$ws3.Rows(12).RowHeight = 28.5

# Row 13 starts a new (still mostly empty) B787 row -- also A13 goes A321 -> B787
$ws3.Range("A13").Value = "B787"
$ws3.Range("B13").Value = "Y SEAT"

# Rows 14-17 were A321, re-tagged to ATR72 (contents otherwise untouched)
$ws3.Range("A14").Value = "ATR72"
$ws3.Range("A15").Value = "ATR72"
$ws3.Range("A16").Value = "ATR72"
$ws3.Range("A17").Value = "ATR72"

# Six new blank A321 rows appended, copying row 17's formatting
$ws3.Range("A17:F17").Copy()
$ws3.Range("A18:F23").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("A18").Value = "A321"
$ws3.Range("A19").Value = "A321"
$ws3.Range("A20").Value = "A321"
$ws3.Range("A21").Value = "A321"
$ws3.Range("A22").Value = "A321"
$ws3.Range("A23").Value = "A321"
$excel.CutCopyMode = $false

# Pax seat becomes the selected/active tab, with F19 highlighted
$ws3.Activate()
$ws3.Range("F19").Select()
